$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Update the letter date: "September 19, 2025" -> "September 21, 2025"
# ---------------------------------------------------------------------
$dateRng = $d.Content
$dateRng.Find.ClearFormatting()
$found1 = $dateRng.Find.Execute("September 19, 2025", $false, $false, $false, $false, `
    $false, $true, 1, $false, "September 21, 2025", 2)
Write-Host "date replace found/executed: $found1"

# ---------------------------------------------------------------------
# 2. Split the recipient mailing-address paragraph
#    "175 Lewis Road Suite, San Jose CA 95111"
#    into two paragraphs:
#       "175 Lewis Road Suite"
#       "San Jose, CA 95111"
#    (only the first occurrence - the recipient block right after
#    "Girma Bekele" - not the PROPERTY ADDRESS table entry further down
#    in the document, which must stay untouched)
# ---------------------------------------------------------------------
$addrRng = $d.Content
$addrRng.Find.ClearFormatting()
$found2 = $addrRng.Find.Execute("175 Lewis Road Suite, San Jose CA 95111", $false, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $addrRng.Text = "175 Lewis Road Suite" + [char]13 + "San Jose, CA 95111"
}
Write-Host "address split found/executed: $found2"

# ---------------------------------------------------------------------
# 3. Remove the blank "NoSpacing" paragraph that used to immediately
#    follow "...Board of Directors"
# ---------------------------------------------------------------------
$bodRng = $d.Content
$bodRng.Find.ClearFormatting()
$found3 = $bodRng.Find.Execute("Board of Directors", $false, $false, $false, $false, `
    $false, $true, 1, $false, "", 0)

$targetIndex = -1
if ($found3) {
    $foundEnd = $bodRng.End
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $pp = $d.Paragraphs.Item($i)
        if ($pp.Range.Start -le $foundEnd -and $pp.Range.End -ge $foundEnd) {
            $targetIndex = $i
            break
        }
    }
}
if ($targetIndex -gt 0) {
    $nextPara = $d.Paragraphs.Item($targetIndex + 1)
    if ($nextPara.Range.Text -eq "" -or $nextPara.Range.Text -eq [char]13) {
        $nextPara.Range.Delete()
    }
}
Write-Host "blank paragraph after Board of Directors removed (index $targetIndex)"

Write-Host "done"
